$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(192,193),
    @(219,220),
    @(227,228),
    @(232,233),
    @(243,244),
    @(364,365),
    @(382,383),
    @(442,443),
    @(463,464),
    @(473,474),
    @(572,573)
)

foreach ($pair in $pairs) {
    $a = $pair[0]
    $b = $pair[1]

    foreach ($col in @("B","E","F","G")) {
        $cellA = $ws.Range("$col$a")
        $cellB = $ws.Range("$col$b")
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}
